$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "urn:nasa:pds:lab_shocked_feldspars::11.0"
$ws.Range("G3").Value = "urn:nasa:pds:lab_shocked_feldspars_21::1.0"
$ws.Range("G4").Value = "urn:nasa:pds:lab_shocked_feldspars_31::1.0"

$ws.Range("G5").Select()
